$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text storage; Excel would
# otherwise auto-coerce plain decimal-looking strings (e.g. "212.87") into
# numbers. Apply a Text format to the whole data range first, write the
# values, then restore the Normal style so no stray per-cell formatting is
# left behind (matches the un-styled cells in the source file).
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '27.673.51'
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').Value = '1.639.33'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '212.87'
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('D6').Value = '0.524'
$ws.Range('E6').Value = '  -0.63%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '23.10'
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  +0.43%  '
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('D11').Value = '0.0895'
$ws.Range('E11').Value = '  +0.66%  '
$ws.Range('D12').Value = '1.872.28'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').Value = '1.660.91'
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('D14').Value = '4.04'
$ws.Range('E14').Value = '  +0.12%  '
$ws.Range('E15').Value = '  -1.72%  '
$ws.Range('D16').Value = '64.62'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').Value = '27.649.55'
$ws.Range('D18').Value = '230.43'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').Value = '7.72'
$ws.Range('E19').Value = '  +1.86%  '
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('E22').Value = '  -1.10%  '
$ws.Range('E23').Value = '  +3.67%  '
$ws.Range('E24').Value = '  -2.87%  '
$ws.Range('D25').Value = '149.67'
$ws.Range('E25').Value = '  +1.94%  '
$ws.Range('D26').Value = '6.93'
$ws.Range('E26').Value = '  -0.93%  '
$ws.Range('E27').Value = '  -1.08%  '
$ws.Range('E28').Value = '  +0.63%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  +0.40%  '
$ws.Range('E31').Value = '  +0.42%  '
$ws.Range('E32').Value = '  +0.50%  '
$ws.Range('D33').Value = '1.445.40'
$ws.Range('E34').Value = '  -0.74%  '
$ws.Range('E35').Value = '  -0.33%  '
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('E37').Value = '  +0.45%  '
$ws.Range('E38').Value = '  -0.63%  '
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('D40').Value = '0.899'
$ws.Range('E40').Value = '  +12.56%  '
$ws.Range('D41').Value = '70.27'
$ws.Range('E41').Value = '  +8.98%  '
$ws.Range('D42').Value = '1.02'
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('E44').Value = '  +1.71%  '
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').Value = '1.781.99'
$ws.Range('E48').Value = '  +3.56%  '
$ws.Range('D49').Value = '86.18'
$ws.Range('E49').Value = '  -1.91%  '
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('E51').Value = '  -0.03%  '

$ws.Range('D2:D51').Style = 'Normal'

